$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures per commit "Updated cryptos list".
# D-column price cells are forced to Text format before assignment so that
# values such as "1.001" or "248.18" are not auto-coerced into numbers by
# Excel's smart entry, then the cell style is reset back to Normal so no
# stray formatting is left behind (matches original inlineStr cells, which
# carried no explicit style).
$rowData = @(
    @{ Row = 2; D = '30.476.60'; E = '  +0.29%  ' }
    @{ Row = 3; D = '1.871.39'; E = '  -0.51%  ' }
    @{ Row = 4; D = $null; E = '  +0.20%  ' }
    @{ Row = 5; D = '248.18'; E = '  +1.09%  ' }
    @{ Row = 6; D = $null; E = '  +0.21%  ' }
    @{ Row = 7; D = '0.4719'; E = '  -0.41%  ' }
    @{ Row = 8; D = '0.2899'; E = '  +0.11%  ' }
    @{ Row = 9; D = '0.06456'; E = '  -1.30%  ' }
    @{ Row = 10; D = '22.04'; E = '  +3.85%  ' }
    @{ Row = 11; D = '0.07696'; E = $null }
    @{ Row = 12; D = '0.7391'; E = '  +0.61%  ' }
    @{ Row = 13; D = '95.92'; E = '  +0.13%  ' }
    @{ Row = 14; D = '1.869.69'; E = '  -0.50%  ' }
    @{ Row = 15; D = '5.155'; E = '  +0.16%  ' }
    @{ Row = 16; D = '272.64'; E = '  -1.43%  ' }
    @{ Row = 17; D = '30.519.88'; E = '  +0.49%  ' }
    @{ Row = 18; D = '13.25'; E = '  -1.22%  ' }
    @{ Row = 19; D = '1.001'; E = '  +0.07%  ' }
    @{ Row = 20; D = '0.000007462'; E = '  -1.70%  ' }
    @{ Row = 21; D = '2.119.19'; E = '  +0.08%  ' }
    @{ Row = 22; D = $null; E = '  +0.27%  ' }
    @{ Row = 23; D = '5.242'; E = '  -0.53%  ' }
    @{ Row = 24; D = '6.163'; E = '  -0.37%  ' }
    @{ Row = 25; D = '9.176'; E = '  -1.05%  ' }
    @{ Row = 26; D = '164.29'; E = '  -0.43%  ' }
    @{ Row = 27; D = '18.68'; E = '  -1.69%  ' }
    @{ Row = 28; D = '1.899'; E = '  -2.16%  ' }
    @{ Row = 29; D = '0.09973'; E = '  +0.64%  ' }
    @{ Row = 30; D = '1.347'; E = '  -2.90%  ' }
    @{ Row = 31; D = '1.512'; E = '  -0.54%  ' }
    @{ Row = 32; D = '4.238'; E = '  -2.20%  ' }
    @{ Row = 33; D = '4.078'; E = '  +0.25%  ' }
    @{ Row = 34; D = '0.04779'; E = '  +0.09%  ' }
    @{ Row = 35; D = '1.116'; E = '  -1.12%  ' }
    @{ Row = 36; D = '0.6912'; E = '  -1.38%  ' }
    @{ Row = 37; D = '2.721'; E = '  +0.17%  ' }
    @{ Row = 38; D = '0.01849'; E = '  -0.12%  ' }
    @{ Row = 39; D = '2.749'; E = '  -0.25%  ' }
    @{ Row = 40; D = '6.231'; E = '  -3.26%  ' }
    @{ Row = 41; D = '73.10'; E = '  +3.97%  ' }
    @{ Row = 42; D = '1.964'; E = '  +2.23%  ' }
    @{ Row = 43; D = $null; E = '  +0.19%  ' }
    @{ Row = 44; D = '0.4147'; E = '  -0.52%  ' }
    @{ Row = 45; D = '0.8331'; E = '  -1.38%  ' }
    @{ Row = 46; D = '101.21'; E = '  -1.31%  ' }
    @{ Row = 47; D = '9.351'; E = '  -0.32%  ' }
    @{ Row = 48; D = '35.30'; E = '  -0.22%  ' }
    @{ Row = 49; D = '6.967'; E = '  -2.13%  ' }
    @{ Row = 50; D = '911.77'; E = '  -2.02%  ' }
    @{ Row = 51; D = '0.05654'; E = '  +1.27%  ' }
)

foreach ($item in $rowData) {
    if ($null -ne $item.D) {
        $cell = $ws.Cells.Item($item.Row, 4)   # column D
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.Style = "Normal"
    }
    if ($null -ne $item.E) {
        $ws.Cells.Item($item.Row, 5).Value = $item.E   # column E
    }
}
